$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting (no numeric auto-conversion)
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '19.955.92'
$ws.Range("E2").Value = '  -8.10%  '
$ws.Range("D3").Value = '1.405.93'
$ws.Range("E3").Value = '  -8.65%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").Value = '272.52'
$ws.Range("E6").Value = '  -5.96%  '
$ws.Range("D7").Value = '0.3677'
$ws.Range("E7").Value = '  -6.91%  '
$ws.Range("D8").Value = '0.3128'
$ws.Range("E8").Value = '  -2.35%  '
$ws.Range("D9").Value = '39.73'
$ws.Range("E9").Value = '  -6.87%  '
$ws.Range("D10").Value = '1.012'
$ws.Range("E10").Value = '  -6.61%  '
$ws.Range("D11").Value = '0.06491'
$ws.Range("E11").Value = '  -9.88%  '
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.25%  '
$ws.Range("D13").Value = '5.454'
$ws.Range("E13").Value = '  -5.21%  '
$ws.Range("D14").Value = '17.34'
$ws.Range("E14").Value = '  -6.00%  '
$ws.Range("D15").Value = '6.141'
$ws.Range("E15").Value = '  -7.77%  '
$ws.Range("D16").Value = '1.405.53'
$ws.Range("D17").Value = '0.00001012'
$ws.Range("E17").Value = '  -8.02%  '
$ws.Range("D18").Value = '0.05688'
$ws.Range("E18").Value = '  -13.86%  '
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("D20").Value = '70.76'
$ws.Range("E20").Value = '  -16.24%  '
$ws.Range("D21").Value = '5.564'
$ws.Range("E21").Value = '  -9.68%  '
$ws.Range("E22").Value = '  -5.67%  '
$ws.Range("D23").Value = '11.02'
$ws.Range("E23").Value = '  +1.52%  '
$ws.Range("D24").Value = '2.260'
$ws.Range("E24").Value = '  -4.77%  '
$ws.Range("D25").Value = '19.984.53'
$ws.Range("E25").Value = '  -7.98%  '
$ws.Range("D26").Value = '2.232'
$ws.Range("E26").Value = '  -7.27%  '
$ws.Range("D27").Value = '135.57'
$ws.Range("E27").Value = '  -10.59%  '
$ws.Range("D28").Value = '16.87'
$ws.Range("E28").Value = '  -8.82%  '
$ws.Range("D29").Value = '1.563.83'
$ws.Range("E29").Value = '  -8.83%  '
$ws.Range("D30").Value = '109.23'
$ws.Range("E30").Value = '  -7.36%  '
$ws.Range("D31").Value = '4.112'
$ws.Range("E31").Value = '  -15.27%  '
$ws.Range("D32").Value = '5.296'
$ws.Range("E32").Value = '  -13.57%  '
$ws.Range("D33").Value = '0.8166'
$ws.Range("E33").Value = '  -16.67%  '
$ws.Range("D34").Value = '0.07663'
$ws.Range("E34").Value = '  -5.98%  '
$ws.Range("D35").Value = '8.408'
$ws.Range("E35").Value = '  -2.15%  '
$ws.Range("D36").Value = '1.448'
$ws.Range("E36").Value = '  -2.59%  '
$ws.Range("D37").Value = '0.05791'
$ws.Range("E37").Value = '  -3.58%  '
$ws.Range("D38").Value = '4.827'
$ws.Range("E38").Value = '  -7.50%  '
$ws.Range("D39").Value = '0.9998'
$ws.Range("E39").Value = '  +0.15%  '
$ws.Range("D40").Value = '0.02069'
$ws.Range("E40").Value = '  -7.79%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '0.1900'
$ws.Range("E41").Value = '  -7.47%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").Value = '10.37'
$ws.Range("E42").Value = '  -8.94%  '
$ws.Range("D43").Value = '1.099'
$ws.Range("E43").Value = '  -7.32%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '12.30'
$ws.Range("E44").Value = '  -6.86%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.5274'
$ws.Range("E45").Value = '  -9.82%  '
$ws.Range("D46").Value = '3.508'
$ws.Range("E46").Value = '  -5.93%  '
$ws.Range("D47").Value = '0.5119'
$ws.Range("E47").Value = '  -8.75%  '
$ws.Range("D48").Value = '111.54'
$ws.Range("E48").Value = '  -4.58%  '
$ws.Range("D49").Value = '1.760'
$ws.Range("E49").Value = '  -7.35%  '
$ws.Range("D50").Value = '1.034'
$ws.Range("E50").Value = '  -11.68%  '
$ws.Range("D51").Value = '0.9984'
$ws.Range("E51").Value = '  +0.05%  '

# Restore default (Normal) style so no visual/style differences remain versus the original cells
$priceRange.Style = "Normal"

